{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targets = new Set([\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n]);\n\n// Find the index of the \"LOQ4233\" paragraph so we only remove the blank\n// paragraph that immediately follows it (not every blank paragraph in the\n// document).\nlet anchorIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"LOQ4233: Gest\\u00e3o de Neg\\u00f3cios (Requisito fraco)\") !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nconst toDelete = [];\nif (anchorIndex !== -1 && anchorIndex + 1 < paragraphs.items.length) {\n  // The blank paragraph directly after the LOQ4233 line.\n  if (paragraphs.items[anchorIndex + 1].text.trim() === \"\") {\n    toDelete.push(paragraphs.items[anchorIndex + 1]);\n  }\n}\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (targets.has(paragraphs.items[i].text.trim())) {\n    toDelete.push(paragraphs.items[i]);\n  }\n}\n\nfor (const p of toDelete) {\n  p.delete();\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$targets = @(\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n)\n\n# Walk the paragraphs collection once, recording the ones that must go:\n#   - the two \"footer\" paragraphs identified by their exact text, and\n#   - the blank paragraph that immediately follows the LOQ4233 requirement line.\n$count = $d.Paragraphs.Count\n$toDelete = New-Object System.Collections.ArrayList\n$anchorIndex = -1\n\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.TrimEnd(\"`r\", \"`a\").Trim()\n    if ($t -like \"LOQ4233:*\") {\n        $anchorIndex = $i\n    }\n    if ($targets -contains $t) {\n        [void]$toDelete.Add($i)\n    }\n}\n\nif ($anchorIndex -ge 1 -and ($anchorIndex + 1) -le $count) {\n    $blank = $d.Paragraphs.Item($anchorIndex + 1)\n    $blankText = $blank.Range.Text.TrimEnd(\"`r\", \"`a\").Trim()\n    if ($blankText -eq \"\") {\n        [void]$toDelete.Add($anchorIndex + 1)\n    }\n}\n\n# Delete from the highest index down so earlier indices stay valid.\n$sorted = $toDelete | Sort-Object -Descending -Unique\nforeach ($idx in $sorted) {\n    $d.Paragraphs.Item($idx).Range.Delete()\n}\n"}
